$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new row 5 with ConvMF results (Validation split)
$ws.Range("A5").Value = "ConvMF"
$ws.Range("B5").Value = "Validation"
$ws.Range("C5").Value = 0.73399999999999999
$ws.Range("D5").Value = 0.5615
$ws.Range("E5").Value = 0.003
$ws.Range("F5").Value = 0.0086
$ws.Range("G5").Value = 0.0036
$ws.Range("H5").Value = 0.0062
$ws.Range("I5").Value = 0.0131

# Add the Harmonic mean formulas for rows 4 and 5 (shared formula)
$ws.Range("J4").Formula = "=7/((1/C4)+(1/D4)+(1/E4)+(1/F4)+(1/G4)+(1/H4)+(1/I4))"
$ws.Range("J5").Formula = "=7/((1/C5)+(1/D5)+(1/E5)+(1/F5)+(1/G5)+(1/H5)+(1/I5))"

# Update the active selection to L6, as recorded in the saved workbook view
$ws.Range("L6").Select()

$wb.Save()
